$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'60.102.05"
$ws.Range("E2").Value = "  +1.10%  "

# Row 3
$ws.Range("D3").Value = "'2.590.80"
$ws.Range("E3").Value = "  +0.18%  "

# Row 4
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").Value = "'578.62"
$ws.Range("E5").Value = "  +4.47%  "

# Row 6
$ws.Range("D6").Value = "'142.30"
$ws.Range("E6").Value = "  +1.41%  "

# Row 7
$ws.Range("E7").Value = "  +0.06%  "

# Row 8
$ws.Range("E8").Value = "  +0.40%  "

# Row 9
$ws.Range("D9").Value = "'2.592.26"
$ws.Range("E9").Value = "  -0.34%  "

# Row 10
$ws.Range("E10").Value = "  -2.81%  "

# Row 11
$ws.Range("E11").Value = "  +0.91%  "

# Row 12
$ws.Range("E12").Value = "  -2.48%  "

# Row 13
$ws.Range("D13").Value = "'0.369"
$ws.Range("E13").Value = "  +3.24%  "

# Row 14
$ws.Range("D14").Value = "'3.051.99"
$ws.Range("E14").Value = "  +0.04%  "

# Row 15
$ws.Range("D15").Value = "'24.58"
$ws.Range("E15").Value = "  +6.68%  "

# Row 16
$ws.Range("D16").Value = "'60.092.01"
$ws.Range("E16").Value = "  +1.10%  "

# Row 17
$ws.Range("E17").Value = "  +2.36%  "

# Row 18
$ws.Range("D18").Value = "'2.599.09"
$ws.Range("E18").Value = "  +0.01%  "

# Row 19
$ws.Range("D19").Value = "'11.42"
$ws.Range("E19").Value = "  +9.68%  "

# Row 20
$ws.Range("D20").Value = "'4.62"
$ws.Range("E20").Value = "  +1.32%  "

# Row 21
$ws.Range("D21").Value = "'345.83"
$ws.Range("E21").Value = "  +1.63%  "

# Row 22
$ws.Range("E22").Value = "  +4.41%  "

# Row 23
$ws.Range("E23").Value = "  +0.32%  "

# Row 24
$ws.Range("D24").Value = "'0.524"
$ws.Range("E24").Value = "  +8.31%  "

# Row 25
$ws.Range("D25").Value = "'62.81"
$ws.Range("E25").Value = "  +0.09%  "

# Row 26
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  +0.24%  "

# Row 27
$ws.Range("D27").Value = "'0.158"
$ws.Range("E27").Value = "  +0.13%  "

# Row 28
$ws.Range("D28").Value = "'8.01"
$ws.Range("E28").Value = "  +7.03%  "

# Row 29
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D29").Value = "'0.0₃0787"
$ws.Range("E29").Value = "  +2.43%  "

# Row 30
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'1.85"
$ws.Range("E30").Value = "  +10.24%  "

# Row 31
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  +0.06%  "

# Row 32
$ws.Range("E32").Value = "  +4.05%  "

# Row 33
$ws.Range("D33").Value = "'164.26"
$ws.Range("E33").Value = "  +4.41%  "

# Row 34
$ws.Range("D34").Value = "'19.38"
$ws.Range("E34").Value = "  +0.18%  "

# Row 35
$ws.Range("D35").Value = "'4.26"
$ws.Range("E35").Value = "  +4.42%  "

# Row 36
$ws.Range("D36").Value = "'0.975"
$ws.Range("E36").Value = "  +6.00%  "

# Row 37
$ws.Range("E37").Value = "  +6.83%  "

# Row 38
$ws.Range("E38").Value = "  +9.34%  "

# Row 39
$ws.Range("D39").Value = "'38.01"
$ws.Range("E39").Value = "  +1.15%  "

# Row 40
$ws.Range("D40").Value = "'3.90"
$ws.Range("E40").Value = "  +6.19%  "

# Row 41
$ws.Range("D41").Value = "'308.42"
$ws.Range("E41").Value = "  +6.54%  "

# Row 42
$ws.Range("D42").Value = "'0.837"
$ws.Range("E42").Value = "  -0.17%  "

# Row 43
$ws.Range("D43").Value = "'134.97"
$ws.Range("E43").Value = "  -1.07%  "

# Row 44
$ws.Range("D44").Value = "'0.999"
$ws.Range("E44").Value = "  +0.16%  "

# Row 45
$ws.Range("E45").Value = "  +1.54%  "

# Row 46
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "'4.99"
$ws.Range("E46").Value = "  +10.11%  "

# Row 47
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "'0.601"
$ws.Range("E47").Value = "  +0.38%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'19.63"
$ws.Range("E48").Value = "  +3.76%  "

# Row 49
$ws.Range("D49").Value = "'0.0546"
$ws.Range("E49").Value = "  +2.23%  "

# Row 50
$ws.Range("D50").Value = "'19.97"
$ws.Range("E50").Value = "  +7.08%  "

# Row 51
$ws.Range("E51").Value = "  +2.25%  "

